# Adds a new "VSC Trunking" setting row to the "Common" sheet of the
# deployment workbook, inserted right above the existing "Security"
# section header (i.e. becomes the new row 60; "Security" and everything
# below shifts down by one row). Also adds the corresponding cell
# comment, and re-homes every comment that was anchored on a row at or
# below the old "Security" header so that it follows its cell's content
# to the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

$insertRow = 60
$lastRow = 132

# 1) Shift every comment on rows [insertRow, lastRow] down by one row so
#    it stays attached to the same content after the row insert below.
#    Walk from the bottom up so we never overwrite a comment we still
#    need to read.
for ($r = $lastRow; $r -ge $insertRow; $r--) {
    $src = $ws.Range("A" + $r)
    if ($src.Comment) {
        $txt = $src.Comment.Text()
        $dst = $ws.Range("A" + ($r + 1))
        $dst.AddComment($txt) | Out-Null
        $src.Comment.Delete()
    }
}

# 2) Insert the new blank row, pushing "Security" (and everything after
#    it) down to row 61.
$ws.Rows.Item($insertRow).Insert()

# 3) Give the new row the same look as the other plain setting rows
#    (copy format from the row above, which is a normal label/value pair).
$ws.Range("A" + ($insertRow - 1) + ":B" + ($insertRow - 1)).Copy() | Out-Null
$ws.Range("A" + $insertRow + ":B" + $insertRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 4) Set the label text for the new row; leave the value cell blank.
$ws.Range("A" + $insertRow).Value = "VSC Trunking"
$ws.Range("B" + $insertRow).ClearContents()

# 5) Add the comment describing the new setting.
$comment = "This enables trunking between vsc control port and its underlay ports (openstack_external_port1_name, openstack_external_port2_name,openstack_external_port3_name) as defined under vscs.yml [default: False]"
$ws.Range("A" + $insertRow).AddComment($comment) | Out-Null

# 6) The field is a boolean flag like the other True/False settings, so
#    give its value cell the same data validation dropdown (xlValidateList,
#    xlValidAlertWarning, xlBetween) as its siblings.
$validation = $ws.Range("B" + $insertRow).Validation
$validation.Add(3, 2, 1, """true,false""") | Out-Null
$validation.ErrorTitle = "Invalid Entry"
$validation.ErrorMessage = "Your entry is not true or false, change anyway?"
$validation.InputTitle = "True or False Selection"
$validation.InputMessage = "Please select true or false"
$validation.ShowInput = $true
$validation.ShowError = $true
